$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 34; A = '2024-10-28 11:06:19'; B = 'ORM-0516014'; C = 'MolPort-009-747-996'; D = 'C14H11N3O2S'; E = '/home/robekott/ERAT/MOLPORT/molport_test.sdf'; F = 'Success'; G = 201; H = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f540abc105f6067d39580"},"data":{"type":"material","id":"asset:671f540abc105f6067d39580","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f540abc105f6067d39580"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"671f540abc105f6067d39580","id":"asset:671f540abc105f6067d39580","eid":"asset:671f540abc105f6067d39580","name":"ORM-0516014","synonyms":["N1C(ONC1C)C1CCC(NC(=O)C2SCCC2)CC1","C14H11N3O2S"],"description":"","createdAt":"2024-10-28T09:06:18.567Z","editedAt":"2024-10-28T09:06:18.567Z","type":"asset","digest":"31900071","fields":{"Chemical Name":{"value":"N-[4-(3-methyl-1,2,4-oxadiazol-5-yl)phenyl]thiophene-2-carboxamide"},"Description":{"value":""},"Exact Mass":{"value":"285.0572"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C<sub>14</sub>H<sub>11</sub>N<sub>3</sub>O<sub>2</sub>S"},"Molecular Weight":{"value":"285.32 g/mol"},"Name":{"value":"ORM-0516014"},"Stereochemistry":{"value":"achiral"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:671f540abc105f6067d39581","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:671f540abc105f6067d39581"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f540abc105f6067d39580/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:671f540abc105f6067d39580"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"material","id":"batch:671f540abc105f6067d39581","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:671f540abc105f6067d39581"},"attributes":{"type":"batch","eid":"batch:671f540abc105f6067d39581","name":"ORM-0516014-001","digest":"16551587","fields":{"Batch Chemical Name":{"value":"N-[4-(3-methyl-1,2,4-oxadiazol-5-yl)phenyl]thiophene-2-carboxamide"},"Batch Molecular Formula":{"value":"C<sub>14</sub>H<sub>11</sub>N<sub>3</sub>O<sub>2</sub>S"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"285.32 g/mol"},"Library ID":{"value":""},"Name":{"value":"ORM-0516014-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-10-28T09:06:18.891Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"materialDrawing","id":"asset:671f540abc105f6067d39580","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f540abc105f6067d39580/drawing?format=cdxml"},"attributes":{"id":"asset:671f540abc105f6067d39580","type":"CHEMICAL_DRAWING"}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"15127052","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}}]}' }
    @{ Row = 35; A = '2024-10-28 11:09:41'; B = 'ORM-0516015'; C = 'Z195631098'; D = 'C18H19CL2N5OS'; E = '/home/robekott/ERAT/examples/compound_test.sdf'; F = 'Success'; G = 201; H = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f54d46395b975bca22f75"},"data":{"type":"material","id":"asset:671f54d46395b975bca22f75","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f54d46395b975bca22f75"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"671f54d46395b975bca22f75","id":"asset:671f54d46395b975bca22f75","eid":"asset:671f54d46395b975bca22f75","name":"ORM-0516015","synonyms":["CN(CC(=O)NC1C(CL)CCCC1CL)CC1NC(N)C2C(C)C(C)SC2N1","C18H19CL2N5OS"],"description":"","createdAt":"2024-10-28T09:09:40.547Z","editedAt":"2024-10-28T09:09:40.547Z","type":"asset","digest":"54222663","fields":{"Chemical Name":{"value":"2-(((4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl)methyl)(methyl)amino)-N-(2,6-dichlorophenyl)acetamide"},"Description":{"value":""},"Exact Mass":{"value":"423.06874"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C<sub>18</sub>H<sub>19</sub>Cl<sub>2</sub>N<sub>5</sub>OS"},"Molecular Weight":{"value":"424.34 g/mol"},"Name":{"value":"ORM-0516015"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:671f54d46395b975bca22f76","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:671f54d46395b975bca22f76"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f54d46395b975bca22f75/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:671f54d46395b975bca22f75"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"materialDrawing","id":"asset:671f54d46395b975bca22f75","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f54d46395b975bca22f75/drawing?format=cdxml"},"attributes":{"id":"asset:671f54d46395b975bca22f75","type":"CHEMICAL_DRAWING"}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"87955301","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}},{"type":"material","id":"batch:671f54d46395b975bca22f76","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:671f54d46395b975bca22f76"},"attributes":{"type":"batch","eid":"batch:671f54d46395b975bca22f76","name":"ORM-0516015-001","digest":"75127726","fields":{"Batch Chemical Name":{"value":"2-(((4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl)methyl)(methyl)amino)-N-(2,6-dichlorophenyl)acetamide"},"Batch Molecular Formula":{"value":"C<sub>18</sub>H<sub>19</sub>Cl<sub>2</sub>N<sub>5</sub>OS"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"424.34 g/mol"},"Library ID":{"value":""},"Name":{"value":"ORM-0516015-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-10-28T09:09:40.795Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}}]}' }
    @{ Row = 36; A = '2024-10-28 11:09:44'; B = 'ORM-0516016'; C = 'Z2754556176'; D = 'C17H28N4O2'; E = '/home/robekott/ERAT/examples/compound_test.sdf'; F = 'Success'; G = 201; H = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f54d7bc105f6067d39582"},"data":{"type":"material","id":"asset:671f54d7bc105f6067d39582","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f54d7bc105f6067d39582"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"671f54d7bc105f6067d39582","id":"asset:671f54d7bc105f6067d39582","eid":"asset:671f54d7bc105f6067d39582","name":"ORM-0516016","synonyms":["COCCN1CCNC1C1N(CC1(C)C)C(=O)CC1(CN)CC1","C17H28N4O2"],"description":"","createdAt":"2024-10-28T09:09:43.906Z","editedAt":"2024-10-28T09:09:43.906Z","type":"asset","digest":"33035143","fields":{"Chemical Name":{"value":"2-(1-(aminomethyl)cyclopropyl)-1-(2-(1-(2-methoxyethyl)-1H-imidazol-2-yl)-3,3-dimethylazetidin-1-yl)ethan-1-one"},"Description":{"value":""},"Exact Mass":{"value":"320.22123"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C<sub>17</sub>H<sub>28</sub>N<sub>4</sub>O<sub>2</sub>"},"Molecular Weight":{"value":"320.44 g/mol"},"Name":{"value":"ORM-0516016"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:671f54d8bc105f6067d39583","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:671f54d8bc105f6067d39583"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f54d7bc105f6067d39582/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:671f54d7bc105f6067d39582"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"material","id":"batch:671f54d8bc105f6067d39583","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:671f54d8bc105f6067d39583"},"attributes":{"type":"batch","eid":"batch:671f54d8bc105f6067d39583","name":"ORM-0516016-001","digest":"13092146","fields":{"Batch Chemical Name":{"value":"2-(1-(aminomethyl)cyclopropyl)-1-(2-(1-(2-methoxyethyl)-1H-imidazol-2-yl)-3,3-dimethylazetidin-1-yl)ethan-1-one, hydrogen bromide"},"Batch Molecular Formula":{"value":"C<sub>17</sub>H<sub>28</sub>N<sub>4</sub>O<sub>2</sub>&middot;BRH"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"359.901 g/mol"},"Library ID":{"value":""},"Name":{"value":"ORM-0516016-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-10-28T09:09:44.125Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"materialDrawing","id":"asset:671f54d7bc105f6067d39582","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f54d7bc105f6067d39582/drawing?format=cdxml"},"attributes":{"id":"asset:671f54d7bc105f6067d39582","type":"CHEMICAL_DRAWING"}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"87955301","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}}]}' }
    @{ Row = 37; A = '2024-10-28 11:11:08'; B = 'ORM-0516017'; C = 'Z195631098'; D = 'C18H19CL2N5OS'; E = '/home/robekott/ERAT/examples/compound_test.sdf'; F = 'Success'; G = 201; H = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f552b6395b975bca22f77"},"data":{"type":"material","id":"asset:671f552b6395b975bca22f77","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f552b6395b975bca22f77"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"671f552b6395b975bca22f77","id":"asset:671f552b6395b975bca22f77","eid":"asset:671f552b6395b975bca22f77","name":"ORM-0516017","synonyms":["CN(CC(=O)NC1C(CL)CCCC1CL)CC1NC(N)C2C(C)C(C)SC2N1","C18H19CL2N5OS"],"description":"","createdAt":"2024-10-28T09:11:07.783Z","editedAt":"2024-10-28T09:11:07.783Z","type":"asset","digest":"88092569","fields":{"Chemical Name":{"value":"2-(((4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl)methyl)(methyl)amino)-N-(2,6-dichlorophenyl)acetamide"},"Description":{"value":""},"Exact Mass":{"value":"423.06874"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C<sub>18</sub>H<sub>19</sub>Cl<sub>2</sub>N<sub>5</sub>OS"},"Molecular Weight":{"value":"424.34 g/mol"},"Name":{"value":"ORM-0516017"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:671f552b6395b975bca22f78","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:671f552b6395b975bca22f78"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f552b6395b975bca22f77/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:671f552b6395b975bca22f77"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"materialDrawing","id":"asset:671f552b6395b975bca22f77","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f552b6395b975bca22f77/drawing?format=cdxml"},"attributes":{"id":"asset:671f552b6395b975bca22f77","type":"CHEMICAL_DRAWING"}},{"type":"material","id":"batch:671f552b6395b975bca22f78","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:671f552b6395b975bca22f78"},"attributes":{"type":"batch","eid":"batch:671f552b6395b975bca22f78","name":"ORM-0516017-001","digest":"64567284","fields":{"Batch Chemical Name":{"value":"2-(((4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl)methyl)(methyl)amino)-N-(2,6-dichlorophenyl)acetamide"},"Batch Molecular Formula":{"value":"C<sub>18</sub>H<sub>19</sub>Cl<sub>2</sub>N<sub>5</sub>OS"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"424.34 g/mol"},"Library ID":{"value":""},"Name":{"value":"ORM-0516017-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-10-28T09:11:08.029Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"68126978","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}}]}' }
    @{ Row = 38; A = '2024-10-28 11:11:11'; B = 'ORM-0516018'; C = 'Z2754556176'; D = 'C17H28N4O2'; E = '/home/robekott/ERAT/examples/compound_test.sdf'; F = 'Success'; G = 201; H = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f552e1a12207a92bcb007"},"data":{"type":"material","id":"asset:671f552e1a12207a92bcb007","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f552e1a12207a92bcb007"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"671f552e1a12207a92bcb007","id":"asset:671f552e1a12207a92bcb007","eid":"asset:671f552e1a12207a92bcb007","name":"ORM-0516018","synonyms":["COCCN1CCNC1C1N(CC1(C)C)C(=O)CC1(CN)CC1","C17H28N4O2"],"description":"","createdAt":"2024-10-28T09:11:10.782Z","editedAt":"2024-10-28T09:11:10.782Z","type":"asset","digest":"41788883","fields":{"Chemical Name":{"value":"2-(1-(aminomethyl)cyclopropyl)-1-(2-(1-(2-methoxyethyl)-1H-imidazol-2-yl)-3,3-dimethylazetidin-1-yl)ethan-1-one"},"Description":{"value":""},"Exact Mass":{"value":"320.22123"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C<sub>17</sub>H<sub>28</sub>N<sub>4</sub>O<sub>2</sub>"},"Molecular Weight":{"value":"320.44 g/mol"},"Name":{"value":"ORM-0516018"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:671f552e1a12207a92bcb008","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:671f552e1a12207a92bcb008"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f552e1a12207a92bcb007/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:671f552e1a12207a92bcb007"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"material","id":"batch:671f552e1a12207a92bcb008","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:671f552e1a12207a92bcb008"},"attributes":{"type":"batch","eid":"batch:671f552e1a12207a92bcb008","name":"ORM-0516018-001","digest":"58492800","fields":{"Batch Chemical Name":{"value":"2-(1-(aminomethyl)cyclopropyl)-1-(2-(1-(2-methoxyethyl)-1H-imidazol-2-yl)-3,3-dimethylazetidin-1-yl)ethan-1-one, hydrogen bromide"},"Batch Molecular Formula":{"value":"C<sub>17</sub>H<sub>28</sub>N<sub>4</sub>O<sub>2</sub>&middot;BRH"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"359.901 g/mol"},"Library ID":{"value":""},"Name":{"value":"ORM-0516018-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-10-28T09:11:11.028Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"materialDrawing","id":"asset:671f552e1a12207a92bcb007","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f552e1a12207a92bcb007/drawing?format=cdxml"},"attributes":{"id":"asset:671f552e1a12207a92bcb007","type":"CHEMICAL_DRAWING"}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"68126978","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}}]}' }
)

foreach ($rowData in $rows) {
    $r = $rowData.Row
    $ws.Cells.Item($r, 1).Value = $rowData.A
    $ws.Cells.Item($r, 2).Value = $rowData.B
    $ws.Cells.Item($r, 3).Value = $rowData.C
    $ws.Cells.Item($r, 4).Value = $rowData.D
    $ws.Cells.Item($r, 5).Value = $rowData.E
    $ws.Cells.Item($r, 6).Value = $rowData.F
    $ws.Cells.Item($r, 7).Value = $rowData.G
    $ws.Cells.Item($r, 8).Value = $rowData.H
}

Write-Host "Added rows 34-38"
